$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D (Price) to retain text formatting so that numeric-looking
# strings (e.g. "1.001", "0.7684") are not auto-converted to numbers by Excel.
$ws.Range("D2:D51").NumberFormat = "@"

# Row 51: coin changed from Cronos to Maker
$ws.Range("B51").Value = 'Maker'
$ws.Range("C51").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'

$ws.Range("D2").Value = '29.827.39'
$ws.Range("E2").Value = '  -1.04%  '
$ws.Range("D3").Value = '1.900.24'
$ws.Range("E3").Value = '  -0.60%  '
$ws.Range("D4").Value = '1.001'
$ws.Range("E4").Value = '  +0.03%  '
$ws.Range("D5").Value = '0.7684'
$ws.Range("E5").Value = '  +3.82%  '
$ws.Range("D6").Value = '240.31'
$ws.Range("E6").Value = '  -1.61%  '
$ws.Range("D7").Value = '1.001'
$ws.Range("E7").Value = '  +0.04%  '
$ws.Range("D8").Value = '0.3054'
$ws.Range("E8").Value = '  -2.43%  '
$ws.Range("D9").Value = '25.45'
$ws.Range("E9").Value = '  -5.27%  '
$ws.Range("D10").Value = '0.06848'
$ws.Range("E10").Value = '  -1.92%  '
$ws.Range("D11").Value = '0.07984'
$ws.Range("E11").Value = '  -0.14%  '
$ws.Range("D12").Value = '1.905.68'
$ws.Range("E12").Value = '  +0.02%  '
$ws.Range("D13").Value = '0.7371'
$ws.Range("E13").Value = '  -5.62%  '
$ws.Range("D14").Value = '5.174'
$ws.Range("E14").Value = '  -2.15%  '
$ws.Range("D15").Value = '91.35'
$ws.Range("E15").Value = '  -1.00%  '
$ws.Range("D16").Value = '29.852.01'
$ws.Range("E16").Value = '  -0.96%  '
$ws.Range("D17").Value = '13.74'
$ws.Range("E17").Value = '  -4.62%  '
$ws.Range("D18").Value = '5.888'
$ws.Range("E18").Value = '  -0.67%  '
$ws.Range("D19").Value = '246.77'
$ws.Range("E19").Value = '  +1.99%  '
$ws.Range("D20").Value = '0.000007709'
$ws.Range("E20").Value = '  -1.77%  '
$ws.Range("D21").Value = '1.001'
$ws.Range("E21").Value = '  -0.01%  '
$ws.Range("D22").Value = '2.149.98'
$ws.Range("E22").Value = '  -0.29%  '
$ws.Range("D23").Value = '1.001'
$ws.Range("E23").Value = '  +0.13%  '
$ws.Range("D24").Value = '6.952'
$ws.Range("E24").Value = '  -4.00%  '
$ws.Range("D25").Value = '166.72'
$ws.Range("E25").Value = '  -0.91%  '
$ws.Range("D26").Value = '9.268'
$ws.Range("E26").Value = '  -2.00%  '
$ws.Range("D27").Value = '18.72'
$ws.Range("E27").Value = '  -2.08%  '
$ws.Range("D28").Value = '0.1284'
$ws.Range("E28").Value = '  -0.06%  '
$ws.Range("D29").Value = '2.029'
$ws.Range("E29").Value = '  -2.16%  '
$ws.Range("D30").Value = '1.396'
$ws.Range("E30").Value = '  +2.71%  '
$ws.Range("D31").Value = '1.512'
$ws.Range("E31").Value = '  -2.38%  '
$ws.Range("D32").Value = '4.272'
$ws.Range("E32").Value = '  -1.86%  '
$ws.Range("D33").Value = '4.061'
$ws.Range("E33").Value = '  -1.20%  '
$ws.Range("D34").Value = '0.05249'
$ws.Range("E34").Value = '  +1.20%  '
$ws.Range("D35").Value = '1.245'
$ws.Range("E35").Value = '  -4.24%  '
$ws.Range("D36").Value = '0.7272'
$ws.Range("E36").Value = '  -3.44%  '
$ws.Range("D37").Value = '2.728'
$ws.Range("E37").Value = '  -0.10%  '
$ws.Range("D38").Value = '0.01909'
$ws.Range("E38").Value = '  -1.85%  '
$ws.Range("D39").Value = '2.777'
$ws.Range("E39").Value = '  -0.61%  '
$ws.Range("D40").Value = '6.195'
$ws.Range("E40").Value = '  -2.87%  '
$ws.Range("D41").Value = '0.4409'
$ws.Range("E41").Value = '  -2.60%  '
$ws.Range("D42").Value = '72.01'
$ws.Range("E42").Value = '  -4.28%  '
$ws.Range("D43").Value = '1.001'
$ws.Range("E43").Value = '  -0.04%  '
$ws.Range("D44").Value = '0.8360'
$ws.Range("E44").Value = '  -0.45%  '
$ws.Range("D45").Value = '1.876'
$ws.Range("E45").Value = '  -4.66%  '
$ws.Range("D46").Value = '7.596'
$ws.Range("E46").Value = '  -3.41%  '
$ws.Range("D47").Value = '100.18'
$ws.Range("E47").Value = '  -1.49%  '
$ws.Range("D48").Value = '9.749'
$ws.Range("E48").Value = '  -1.78%  '
$ws.Range("D49").Value = '2.056.78'
$ws.Range("E49").Value = '  +0.89%  '
$ws.Range("D50").Value = '36.16'
$ws.Range("E50").Value = '  -2.93%  '
$ws.Range("D51").Value = '915.72'
$ws.Range("E51").Value = '  -2.62%  '
